$d = $word.ActiveDocument

# 1) Move the "_GoBack" bookmark from the "Paulo Antonio..." paragraph to the
#    start of the "Area de atividade..." paragraph (right before its runs).
$locate = $d.Content
$locate.Find.ClearFormatting()
$null = $locate.Find.Execute("Área de atividade:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($locate.Start, $locate.Start)

$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $bmRange)

# 2) Merge the " " run and the highlighted "Comida Pronta, LTD." run into a
#    single, unhighlighted run: " Comida Pronta, LTD."
$d.Content.Find.Execute(
    " Comida Pronta, LTD.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Comida Pronta, LTD.",
    2
)

# 3) Merge the runs in the "Garantir que o cliente..." paragraph
$d.Content.Find.Execute(
    "Garantir que o cliente possa efetuar a encomenda online, receber a fatura e poder atualizar os dados. Tem de poder ver as lojas/restaurantes;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Garantir que o cliente possa efetuar a encomenda online, receber a fatura e poder atualizar os dados. Tem de poder ver as lojas/restaurantes;",
    2
)

# 4) Merge the runs in the "Os Entregadores..." paragraph
$d.Content.Find.Execute(
    "Os Entregadores possam ver os locais onde devem ir buscar os produtos e onde devem entregar, bem como confirmar o pagamento.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Os Entregadores possam ver os locais onde devem ir buscar os produtos e onde devem entregar, bem como confirmar o pagamento.",
    2
)
